$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1850.6923
$ws.Range("I17").Value = 1150
$ws.Range("J17").Value = 2162.111
$ws.Range("K17").Value = 3450
$ws.Range("L17").Value = 6486.333
$ws.Range("M17").Value = -3282
$ws.Range("N17").Value = -6822.333
$ws.Range("H80").Value = 1245
$ws.Range("I80").Value = 1285
$ws.Range("K80").Value = 3855
$ws.Range("M80").Value = -2857
$ws.Range("H83").Value = 1245
$ws.Range("I83").Value = 1285
$ws.Range("K83").Value = 11565
$ws.Range("M83").Value = -6573
$ws.Range("H100").Value = 1637.2273
$ws.Range("I100").Value = 1134.4445
$ws.Range("K100").Value = 1134.4445
$ws.Range("M100").Value = -593.4445000000001
$ws.Range("H121").Value = 143591.5
$ws.Range("J121").Value = 154560.08
$ws.Range("L121").Value = 463680.24
$ws.Range("N121").Value = -467174.24
$ws.Range("H129").Value = 2887
$ws.Range("I129").Value = 0
$ws.Range("K129").Value = 0
$ws.Range("M129").ClearContents()
$ws.Range("H131").Value = 3334500
$ws.Range("I131").Value = 3334500
$ws.Range("K131").Value = 10003500
$ws.Range("M131").Value = -9998460
$ws.Range("H134").Value = 60330.75
$ws.Range("J134").Value = 60330.75
$ws.Range("L134").Value = 60330.75
$ws.Range("N134").Value = -70470.75
$ws.Range("H136").Value = 54932.125
$ws.Range("J136").Value = 70114.25
$ws.Range("L136").Value = 70114.25
$ws.Range("N136").Value = -80314.25
$ws.Range("H138").Value = 6586.569
$ws.Range("I138").Value = 4080
$ws.Range("J138").Value = 6859.022
$ws.Range("K138").Value = 12240
$ws.Range("L138").Value = 20577.066
$ws.Range("M138").Value = -7100
$ws.Range("N138").Value = -30857.066
$ws.Range("H141").Value = 4615.067
$ws.Range("I141").Value = 1913.1
$ws.Range("J141").Value = 10019
$ws.Range("K141").Value = 5739.299999999999
$ws.Range("L141").Value = 30057
$ws.Range("M141").Value = -559.2999999999993
$ws.Range("N141").Value = -40417

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H19").Value = 1000
$ws.Range("I19").Value = 1000
$ws.Range("K19").Value = 1000
$ws.Range("M19").Value = -771
$ws.Range("H32").Value = 14483.631
$ws.Range("I32").Value = 8792.514999999999
$ws.Range("J32").Value = 38670.875
$ws.Range("K32").Value = 8792.514999999999
$ws.Range("L32").Value = 38670.875
$ws.Range("M32").Value = -8505.514999999999
$ws.Range("N32").Value = -39244.875
$ws.Range("H41").Value = 2534.25
$ws.Range("I41").Value = 2827.2856
$ws.Range("K41").Value = 2827.2856
$ws.Range("M41").Value = -2413.2856
$ws.Range("H61").Value = 169251.17
$ws.Range("I61").Value = 2498
$ws.Range("K61").Value = 2498
$ws.Range("M61").Value = -2286
$ws.Range("H63").Value = 2542.6
$ws.Range("J63").Value = 5125
$ws.Range("L63").Value = 5125
$ws.Range("N63").Value = -6497
$ws.Range("H66").Value = 2542.6
$ws.Range("J66").Value = 5125
$ws.Range("L66").Value = 25625
$ws.Range("N66").Value = -32489
$ws.Range("H74").Value = 93803.17999999999
$ws.Range("I74").Value = 127510.625
$ws.Range("K74").Value = 127510.625
$ws.Range("M74").Value = -126636.625
$ws.Range("H77").Value = 93803.17999999999
$ws.Range("I77").Value = 127510.625
$ws.Range("K77").Value = 637553.125
$ws.Range("M77").Value = -633185.125
$ws.Range("H97").Value = 633.5294
$ws.Range("I97").Value = 633.5294
$ws.Range("K97").Value = 633.5294
$ws.Range("M97").Value = -137.5294
$ws.Range("H110").Value = 2740
$ws.Range("I110").Value = 2528.5715
$ws.Range("J110").Value = 3233.3333
$ws.Range("K110").Value = 2528.5715
$ws.Range("L110").Value = 3233.3333
$ws.Range("M110").Value = -483.5715
$ws.Range("N110").Value = -7323.3333
$ws.Range("H132").Value = 2964.9285
$ws.Range("I132").Value = 2531.1
$ws.Range("J132").Value = 4049.5
$ws.Range("K132").Value = 7593.299999999999
$ws.Range("L132").Value = 12148.5
$ws.Range("M132").Value = -5063.299999999999
$ws.Range("N132").Value = -17208.5
$ws.Range("H136").Value = 169251.17
$ws.Range("I136").Value = 2498
$ws.Range("K136").Value = 7494
$ws.Range("M136").Value = -4944
$ws.Range("H140").Value = 168330.5
$ws.Range("J140").Value = 168330.5
$ws.Range("L140").Value = 168330.5
$ws.Range("N140").Value = -178690.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 4466293
$ws.Range("I99").Value = 2012.5
$ws.Range("K99").Value = 2012.5
$ws.Range("M99").Value = -514.5
$ws.Range("H107").Value = 5886.077
$ws.Range("I107").Value = 5551.5557
$ws.Range("K107").Value = 5551.5557
$ws.Range("M107").Value = -3631.5557

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1369.4286
$ws.Range("I16").Value = 1002.125
$ws.Range("K16").Value = 1002.125
$ws.Range("M16").Value = -715.125
$ws.Range("H58").Value = 1992.375
$ws.Range("I58").Value = 1848.4286
$ws.Range("J58").Value = 3000
$ws.Range("K58").Value = 1848.4286
$ws.Range("L58").Value = 3000
$ws.Range("M58").Value = -1645.4286
$ws.Range("N58").Value = -3406
$ws.Range("H68").Value = 25268
$ws.Range("I68").Value = 25268
$ws.Range("K68").Value = 25268
$ws.Range("M68").Value = -24519
$ws.Range("H71").Value = 25268
$ws.Range("I71").Value = 25268
$ws.Range("K71").Value = 75804
$ws.Range("M71").Value = -72060
$ws.Range("H93").Value = 42300.375
$ws.Range("I93").Value = 7155.2
$ws.Range("K93").Value = 7155.2
$ws.Range("M93").Value = -5283.2
$ws.Range("H99").Value = 3127472
$ws.Range("I99").Value = 2304
$ws.Range("J99").Value = 6252640
$ws.Range("K99").Value = 2304
$ws.Range("L99").Value = 6252640
$ws.Range("M99").Value = -806
$ws.Range("N99").Value = -6255636
$ws.Range("H113").Value = 1369.4286
$ws.Range("I113").Value = 1002.125
$ws.Range("K113").Value = 1002.125
$ws.Range("M113").Value = 1167.875
$ws.Range("H126").Value = 3127472
$ws.Range("I126").Value = 2304
$ws.Range("J126").Value = 6252640
$ws.Range("K126").Value = 6912
$ws.Range("L126").Value = 18757920
$ws.Range("M126").Value = -4442
$ws.Range("N126").Value = -18762860
$ws.Range("H134").Value = 64018.938
$ws.Range("I134").Value = 1122.5834
$ws.Range("J134").Value = 252708
$ws.Range("K134").Value = 3367.7502
$ws.Range("L134").Value = 758124
$ws.Range("M134").Value = -832.7501999999999
$ws.Range("N134").Value = -763194
$ws.Range("H136").Value = 1992.375
$ws.Range("I136").Value = 1848.4286
$ws.Range("J136").Value = 3000
$ws.Range("K136").Value = 5545.2858
$ws.Range("L136").Value = 9000
$ws.Range("M136").Value = -2995.2858
$ws.Range("N136").Value = -14100
$ws.Range("H141").Value = 451465
$ws.Range("I141").Value = 179000
$ws.Range("K141").Value = 179000
$ws.Range("M141").Value = -173820

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H136").Value = 1245
$ws.Range("I136").Value = 1245
$ws.Range("K136").Value = 3735
$ws.Range("M136").Value = 1365

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 37672.47
$ws.Range("I70").Value = 55498.9
$ws.Range("K70").Value = 55498.9
$ws.Range("M70").Value = -55228.9
$ws.Range("H73").Value = 37672.47
$ws.Range("I73").Value = 55498.9
$ws.Range("K73").Value = 55498.9
$ws.Range("M73").Value = -54562.9
$ws.Range("H122").Value = 10603.315
$ws.Range("I122").Value = 11092.588
$ws.Range("K122").Value = 33277.764
$ws.Range("M122").Value = -30827.764
$ws.Range("H132").Value = 5522
$ws.Range("I132").Value = 5433.76
$ws.Range("J132").Value = 6625
$ws.Range("K132").Value = 16301.28
$ws.Range("L132").Value = 19875
$ws.Range("M132").Value = -13771.28
$ws.Range("N132").Value = -24935

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 2177.375
$ws.Range("I46").Value = 1500
$ws.Range("K46").Value = 1500
$ws.Range("M46").Value = -1312
$ws.Range("H68").Value = 4221.9443
$ws.Range("J68").Value = 3428.5715
$ws.Range("L68").Value = 3428.5715
$ws.Range("N68").Value = -4926.5715
$ws.Range("H71").Value = 4221.9443
$ws.Range("J71").Value = 3428.5715
$ws.Range("L71").Value = 17142.8575
$ws.Range("N71").Value = -24630.8575
$ws.Range("H122").Value = 10718.885
$ws.Range("I122").Value = 12805.059
$ws.Range("J122").Value = 6778.3335
$ws.Range("K122").Value = 38415.177
$ws.Range("L122").Value = 20335.0005
$ws.Range("M122").Value = -35965.177
$ws.Range("N122").Value = -25235.0005
$ws.Range("H132").Value = 4308.7
$ws.Range("I132").Value = 3953.9285
$ws.Range("K132").Value = 11861.7855
$ws.Range("M132").Value = -9331.7855

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 9500
$ws.Range("J62").Value = 9500
$ws.Range("L62").Value = 9500
$ws.Range("N62").Value = -10748
$ws.Range("H65").Value = 9500
$ws.Range("J65").Value = 9500
$ws.Range("L65").Value = 47500
$ws.Range("N65").Value = -53740
$ws.Range("H132").Value = 1212.75
$ws.Range("I132").Value = 1039.3334
$ws.Range("J132").Value = 1733
$ws.Range("K132").Value = 3118.0002
$ws.Range("L132").Value = 5199
$ws.Range("M132").Value = -588.0001999999999
$ws.Range("N132").Value = -10259
$ws.Range("H136").Value = 4546.1035
$ws.Range("I136").Value = 3631.3125
$ws.Range("J136").Value = 5672
$ws.Range("K136").Value = 10893.9375
$ws.Range("L136").Value = 17016
$ws.Range("M136").Value = -8343.9375
$ws.Range("N136").Value = -22116
